$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price and volume(1h) values.
# D-column price cells are forced to text format before assignment so that
# plain-numeric-looking strings (e.g. "1.00", "0.320") are not silently
# converted to numbers by Excel, then the style is reset to Normal so no
# extra cell formatting is introduced.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.814.48'
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.508.51'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.11%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '576.48'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '167.72'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.92%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.507.50'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.20%  '
$ws.Range('E10').Value = '  +3.86%  '
$ws.Range('E11').Value = '  +0.25%  '
$ws.Range('E12').Value = '  +4.15%  '
$ws.Range('E13').Value = '  +2.22%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.952.88'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.35%  '
$ws.Range('E15').Value = '  +3.02%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '69.697.87'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.60%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '24.88'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.73%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.486.06'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.95%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.27'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.35%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.49'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '349.65'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.01%  '
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('E23').Value = '  +0.90%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '70.62'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.56%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.97'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.71%  '
$ws.Range('E27').Value = '  -0.36%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.589.82'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.61%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.35%  '
$ws.Range('E30').Value = '  +0.27%  '
$ws.Range('E31').Value = '  +0.16%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '461.51'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.26%  '
$ws.Range('E33').Value = '  -2.39%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.74'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.19%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.06%  '
$ws.Range('E36').Value = '  +2.30%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '157.53'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.08'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.51'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.89%  '
$ws.Range('E40').Value = '  +0.06%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.320'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.59%  '
$ws.Range('E42').Value = '  -0.21%  '
$ws.Range('E43').Value = '  +1.09%  '
$ws.Range('E45').Value = '  -3.10%  '
$ws.Range('E46').Value = '  -6.22%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '141.48'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.81%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.49'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.06%  '
$ws.Range('E49').Value = '  -0.89%  '
$ws.Range('E50').Value = '  +0.81%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.581'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.38%  '
